# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet that
#    carries a per-language status column (Overview!E:F, zh-cn!C, de-de!C).
# 2) Narrow the two status columns (E/F on Overview, C on zh-cn/de-de) to
#    match the new, shorter header/content width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status values -------------------------------------------------
$overview.Range("E2:E4").Value = "In Translation"
$overview.Range("F2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value     = "In Translation"
$dede.Range("C2:C4").Value     = "In Translation"

# --- Narrow the status columns ---------------------------------------------
# Target stored width is ~13.41 chars; ColumnWidth = 12.5 is the input that
# resolves to the closest width this host's column-width grid can produce.
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = 12.5   # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = 12.5   # column C (Status)
